$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose updated text looks numeric must be forced to Text format
# before assignment, otherwise Excel auto-converts them to numbers and
# mangles values like "1.860.04" or drops trailing zeros like "0.7300".

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.043.17"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.798.57"
$ws.Range("E3").Value = "  -2.44%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("E5").Value = "  -2.05%  "

$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4221"
$ws.Range("E7").Value = "  -2.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3597"
$ws.Range("E8").Value = "  -2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07223"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8436"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.25"
$ws.Range("E11").Value = "  -3.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.04"
$ws.Range("E12").Value = "  -3.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.284"
$ws.Range("E13").Value = "  -3.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.363"
$ws.Range("E14").Value = "  -3.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06819"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.65"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008746"
$ws.Range("E18").Value = "  -3.41%  "

$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.04"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.295.03"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.081"
$ws.Range("E22").Value = "  -0.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.084.01"
$ws.Range("E24").Value = "  -2.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.957"
$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.07"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.18"
$ws.Range("E27").Value = "  -3.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.013"
$ws.Range("E28").Value = "  -5.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.75"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.653"
$ws.Range("E30").Value = "  -11.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08989"
$ws.Range("E31").Value = "  +0.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7300"
$ws.Range("E32").Value = "  -7.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.864"
$ws.Range("E33").Value = "  -3.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.352"
$ws.Range("E34").Value = "  -5.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.093"
$ws.Range("E35").Value = "  -6.90%  "

$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05149"
$ws.Range("E38").Value = "  -5.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01904"
$ws.Range("E39").Value = "  -2.94%  "

$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4975"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.618"
$ws.Range("E42").Value = "  -7.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.087"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.971"
$ws.Range("E44").Value = "  -12.29%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.28"
$ws.Range("E45").Value = "  -3.54%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.16"
$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06319"
$ws.Range("E48").Value = "  -3.52%  "

$ws.Range("E49").Value = "  -5.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.602"
$ws.Range("E50").Value = "  -3.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.726"
$ws.Range("E51").Value = "  -6.49%  "

